$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7 with the refreshed TPM-based NATMI output values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna4"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3747116666666667
$ws.Range("H2").Value = 1.124135
$ws.Range("I2").Value = 0.3914669751594584
$ws.Range("J2").Value = 0.3914669751594584
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.097779
$ws.Range("N2").Value = 78.29333700000001
$ws.Range("O2").Value = 0.9922055808976035
$ws.Range("P2").Value = 0.9922055808976036
$ws.Range("Q2").Value = 9.779142265388336
$ws.Range("R2").Value = 88.01228038849501
$ws.Range("S2").Value = 0.3884157174903181
$ws.Range("T2").Value = 0.3884157174903182

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna4"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3747116666666667
$ws.Range("H3").Value = 1.124135
$ws.Range("I3").Value = 0.3914669751594584
$ws.Range("J3").Value = 0.3914669751594584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.205015
$ws.Range("N3").Value = 0.6150450000000001
$ws.Range("O3").Value = 0.007794419102396499
$ws.Range("P3").Value = 0.007794419102396499
$ws.Range("Q3").Value = 0.07682151234166669
$ws.Range("R3").Value = 0.6913936110750002
$ws.Range("S3").Value = 0.003051257669140258
$ws.Range("T3").Value = 0.003051257669140258

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna4"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.401547
$ws.Range("H4").Value = 1.204641
$ws.Range("I4").Value = 0.4195022558883632
$ws.Range("J4").Value = 0.4195022558883631
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.097779
$ws.Range("N4").Value = 78.29333700000001
$ws.Range("O4").Value = 0.9922055808976035
$ws.Range("P4").Value = 0.9922055808976036
$ws.Range("Q4").Value = 10.479484864113
$ws.Range("R4").Value = 94.31536377701701
$ws.Range("S4").Value = 0.4162324794915684
$ws.Range("T4").Value = 0.4162324794915684

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna4"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.401547
$ws.Range("H5").Value = 1.204641
$ws.Range("I5").Value = 0.4195022558883632
$ws.Range("J5").Value = 0.4195022558883631
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.205015
$ws.Range("N5").Value = 0.6150450000000001
$ws.Range("O5").Value = 0.007794419102396499
$ws.Range("P5").Value = 0.007794419102396499
$ws.Range("Q5").Value = 0.08232315820500002
$ws.Range("R5").Value = 0.7409084238450001
$ws.Range("S5").Value = 0.003269776396794682
$ws.Range("T5").Value = 0.003269776396794682

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efna4"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.18094
$ws.Range("H6").Value = 0.54282
$ws.Range("I6").Value = 0.1890307689521785
$ws.Range("J6").Value = 0.1890307689521785
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.097779
$ws.Range("N6").Value = 78.29333700000001
$ws.Range("O6").Value = 0.9922055808976035
$ws.Range("P6").Value = 0.9922055808976036
$ws.Range("Q6").Value = 4.722132132260001
$ws.Range("R6").Value = 42.49918919034
$ws.Range("S6").Value = 0.1875573839157169
$ws.Range("T6").Value = 0.1875573839157169

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efna4"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.18094
$ws.Range("H7").Value = 0.54282
$ws.Range("I7").Value = 0.1890307689521785
$ws.Range("J7").Value = 0.1890307689521785
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.205015
$ws.Range("N7").Value = 0.6150450000000001
$ws.Range("O7").Value = 0.007794419102396499
$ws.Range("P7").Value = 0.007794419102396499
$ws.Range("Q7").Value = 0.0370954141
$ws.Range("R7").Value = 0.3338587269
$ws.Range("S7").Value = 0.001473385036461559
$ws.Range("T7").Value = 0.001473385036461559

# Remove the now-obsolete rows 8-10 (previously the MuSCs sending-cluster block)
$ws.Range("A8:T10").Delete()

